$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '30.398.48'
$ws.Range('E2').Value = '  -1.29%  '
$ws.Range('D3').Value = '1.892.16'
$ws.Range('E3').Value = '  -1.28%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9996'
$ws.Range('E4').Value = '  -0.26%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '237.84'
$ws.Range('E5').Value = '  -1.51%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.9996'
$ws.Range('E6').Value = '  -0.25%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4844'
$ws.Range('E7').Value = '  -1.44%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2907'
$ws.Range('E8').Value = '  -2.49%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06613'
$ws.Range('E9').Value = '  -2.40%  '
$ws.Range('D10').Value = '1.875.78'
$ws.Range('E10').Value = '  -2.07%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '16.96'
$ws.Range('E11').Value = '  -1.10%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.07337'
$ws.Range('E12').Value = '  -0.53%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.182'
$ws.Range('E13').Value = '  -0.41%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '87.73'
$ws.Range('E14').Value = '  -1.73%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.6633'
$ws.Range('E15').Value = '  -1.66%  '
$ws.Range('D16').Value = '30.344.53'
$ws.Range('E16').Value = '  -1.39%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '13.44'
$ws.Range('E17').Value = '  -1.06%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000007789'
$ws.Range('E18').Value = '  -2.65%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.9998'
$ws.Range('E19').Value = '  -0.22%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '5.446'
$ws.Range('E20').Value = '  +3.17%  '
$ws.Range('D21').Value = '2.138.38'
$ws.Range('E21').Value = '  -0.28%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.9988'
$ws.Range('E22').Value = '  -0.35%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '193.56'
$ws.Range('E23').Value = '  -4.89%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '6.193'
$ws.Range('E24').Value = '  -1.67%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '9.403'
$ws.Range('E25').Value = '  -2.75%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '164.51'
$ws.Range('E26').Value = '  +1.87%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '18.22'
$ws.Range('E27').Value = '  -3.50%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.942'
$ws.Range('E28').Value = '  -1.72%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.446'
$ws.Range('E29').Value = '  +0.64%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.319'
$ws.Range('E30').Value = '  -1.05%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.09145'
$ws.Range('E31').Value = '  -0.80%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.059'
$ws.Range('E32').Value = '  -0.72%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.05097'
$ws.Range('E33').Value = '  -6.54%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.7318'
$ws.Range('E34').Value = '  -2.71%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.148'
$ws.Range('E35').Value = '  +2.27%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.703'
$ws.Range('E36').Value = '  +0.03%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.01795'
$ws.Range('E37').Value = '  -3.93%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.651'
$ws.Range('E38').Value = '  -2.85%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.9186'
$ws.Range('E39').Value = '  -1.20%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.075'
$ws.Range('E40').Value = '  -0.40%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '5.899'
$ws.Range('E41').Value = '  -0.69%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '105.94'
$ws.Range('E42').Value = '  -1.72%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.4319'
$ws.Range('E43').Value = '  -4.23%  '
$ws.Range('E44').Value = '  -0.01%  '
$ws.Range('E45').Value = '  -3.34%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.1332'
$ws.Range('E46').Value = '  -4.36%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.588'
$ws.Range('E47').Value = '  +9.77%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '65.24'
$ws.Range('E48').Value = '  -10.69%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '8.926'
$ws.Range('E49').Value = '  -2.28%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '34.01'
$ws.Range('E50').Value = '  -5.63%  '
$ws.Range('E51').Value = '  -3.77%  '
